$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Cryptocurrency price/volume update (GitHub Actions data refresh).
# D and E columns are stored as text in the source sheet, so force the
# "@" text number format before writing to avoid Excel auto-coercing the
# numeric-looking / percentage-looking strings into real numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "329.46"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.32%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "44.19"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.38%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.503"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.32%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08133"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.58%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.070"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "8.32%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9627"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.17%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.1130"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-4.54%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1885"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.70%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "10.13"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.29%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09972"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.16%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04665"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "4.24%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1059"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.77%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001250"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-2.16%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.04090"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-2.71%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005878"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.54%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.377"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.23%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.408"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2.39%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.669"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "4.46%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3311"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-4.86%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1385"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-2.19%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001299"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "4.39%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004369"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.51%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001244"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "4.51%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003718"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-6.56%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02682"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-0.59%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05656"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "1.82%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007590"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.18%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1409"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.08%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007319"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-8.33%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001975"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-2.00%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008244"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-1.80%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00007035"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-1.61%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000746"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.47%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0005768"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.74%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002505"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "10.47%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003326"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-14.57%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002088"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.47%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001988"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.47%"
